# [Kadastro App] Yeni kayit eklendi: 3019
# Adds a new record (row 78) to both the "Kayitlar" summary sheet and the
# "Erdemli" unit sheet, mirroring the existing layout (Kayit No, Tarih,
# Birim, Parsel Sayisi, Is, Personeller). All values are stored as text,
# matching every other row already on these sheets.

$wb = $excel.ActiveWorkbook

$newRow = 78
$values = @(
    "3019",
    "2025-09-11",
    "Erdemli",
    "1",
    "3B",
    "EMİNE ALANLI KIRCILI (K.Mühendisi), SERDAR ARSLAN (Tekniker)"
)

$sheetNames = @("Kayitlar", "Erdemli")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($col = 1; $col -le $values.Length; $col++) {
        $cell = $ws.Cells.Item($newRow, $col)
        # Force text storage (so numeric-looking values like "3019"/"1" and
        # the date-looking "2025-09-11" stay literal strings, not numbers
        # or dates), then drop back to the default "Normal" style so no
        # extra cell formatting is introduced.
        $cell.NumberFormat = "@"
        $cell.Value = $values[$col - 1]
        $cell.Style = "Normal"
    }
}
